$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; existing rows 32..124 shift down to 33..125
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record
$ws.Cells.Item(32, 1).Value = 3
$ws.Cells.Item(32, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44560
$ws.Cells.Item(32, 5).Value = 5
$ws.Cells.Item(32, 6).Value = 100112026
$ws.Cells.Item(32, 7).Value = "Haba"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 110
$ws.Cells.Item(32, 11).Value = 7500
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 7773
$ws.Cells.Item(32, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 311
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
